# Generate Report for Archive
#
# 1) Status text "Ready for handoff" -> "In Translation" on every sheet that
#    surfaces it (Overview!E2:F2, and the per-language "Status" column on the
#    zh-cn / de-de detail sheets — all four cells share the same string).
# 2) Narrow the "Status" column(s) that held that text: Overview!E:F and the
#    "Status" column (C) on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

# --- 1) Update the status text -------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- 2) Narrow the matching columns ---------------------------------------
# Excel's ColumnWidth is quantized to whole pixels, so the narrowest
# achievable width nearest the target (~13.41 chars) is 13.33 chars; a
# ColumnWidth input of 12.5 lands reliably in that bucket.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
